$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3531.2
$ws.Range("I34").Value = 1701.3334
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 1701.3334
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -1498.3334
$ws.Range("N34").Value = -20406

$ws.Range("H36").Value = 3531.2
$ws.Range("I36").Value = 1701.3334
$ws.Range("J36").Value = 20000
$ws.Range("K36").Value = 1701.3334
$ws.Range("L36").Value = 20000
$ws.Range("M36").Value = -986.3334
$ws.Range("N36").Value = -21430

$ws.Range("H64").Value = 3379.6
$ws.Range("J64").Value = 4611.1113
$ws.Range("L64").Value = 4611.1113
$ws.Range("N64").Value = -5107.1113

$ws.Range("H67").Value = 3379.6
$ws.Range("J67").Value = 4611.1113
$ws.Range("L67").Value = 4611.1113
$ws.Range("N67").Value = -6327.1113

$ws.Range("H74").Value = 4496.25
$ws.Range("I74").Value = 4490
$ws.Range("J74").Value = 4500
$ws.Range("K74").Value = 4490
$ws.Range("L74").Value = 4500
$ws.Range("M74").Value = -3554
$ws.Range("N74").Value = -6372

$ws.Range("H77").Value = 4496.25
$ws.Range("I77").Value = 4490
$ws.Range("J77").Value = 4500
$ws.Range("K77").Value = 22450
$ws.Range("L77").Value = 22500
$ws.Range("M77").Value = -17770
$ws.Range("N77").Value = -31860

$ws.Range("H137").Value = 1226.0869
$ws.Range("I137").Value = 1307.0588
$ws.Range("J137").Value = 996.6667
$ws.Range("K137").Value = 3921.1764
$ws.Range("L137").Value = 2990.0001
$ws.Range("M137").Value = -1371.1764
$ws.Range("N137").Value = -8090.0001


# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2129.7058
$ws.Range("I45").Value = 1322.6
$ws.Range("J45").Value = 3282.7144
$ws.Range("K45").Value = 1322.6
$ws.Range("L45").Value = 3282.7144
$ws.Range("M45").Value = -945.5999999999999
$ws.Range("N45").Value = -4036.7144

$ws.Range("H88").Value = 1991.4
$ws.Range("I88").Value = 1500
$ws.Range("J88").Value = 2114.25
$ws.Range("K88").Value = 1500
$ws.Range("L88").Value = 2114.25
$ws.Range("M88").Value = -1094
$ws.Range("N88").Value = -2926.25

$ws.Range("H91").Value = 1991.4
$ws.Range("I91").Value = 1500
$ws.Range("J91").Value = 2114.25
$ws.Range("K91").Value = 1500
$ws.Range("L91").Value = 2114.25
$ws.Range("M91").Value = -96
$ws.Range("N91").Value = -4922.25

$ws.Range("H122").Value = 2049.4546
$ws.Range("I122").Value = 1667.6923
$ws.Range("K122").Value = 5003.0769
$ws.Range("M122").Value = -2553.0769


# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3253.2856
$ws.Range("I86").Value = 3649.4285
$ws.Range("J86").Value = 2857.1428
$ws.Range("K86").Value = 3649.4285
$ws.Range("L86").Value = 2857.1428
$ws.Range("M86").Value = -2526.4285
$ws.Range("N86").Value = -5103.1428

$ws.Range("H89").Value = 3253.2856
$ws.Range("I89").Value = 3649.4285
$ws.Range("J89").Value = 2857.1428
$ws.Range("K89").Value = 18247.1425
$ws.Range("L89").Value = 14285.714
$ws.Range("M89").Value = -12631.1425
$ws.Range("N89").Value = -25517.714

$ws.Range("H99").Value = 1177.5333
$ws.Range("I99").Value = 1196.6666
$ws.Range("K99").Value = 1196.6666
$ws.Range("M99").Value = 301.3334


# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H36").Value = 20100.25
$ws.Range("I36").Value = 3048
$ws.Range("J36").Value = 25784.334
$ws.Range("K36").Value = 3048
$ws.Range("L36").Value = 25784.334
$ws.Range("M36").Value = -2660
$ws.Range("N36").Value = -26560.334

$ws.Range("H40").Value = 20100.25
$ws.Range("I40").Value = 3048
$ws.Range("J40").Value = 25784.334
$ws.Range("K40").Value = 3048
$ws.Range("L40").Value = 25784.334
$ws.Range("M40").Value = -2888
$ws.Range("N40").Value = -26104.334

$ws.Range("H138").Value = 40592.418
$ws.Range("J138").Value = 40592.418
$ws.Range("L138").Value = 40592.418
$ws.Range("N138").Value = -50872.418

$ws.Range("H140").Value = 90126
$ws.Range("J140").Value = 90126
$ws.Range("L140").Value = 90126
$ws.Range("N140").Value = -100486


# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 712.5714
$ws.Range("I113").Value = 666.2222
$ws.Range("J113").Value = 796
$ws.Range("K113").Value = 1998.6666
$ws.Range("L113").Value = 2388
$ws.Range("M113").Value = 171.3334
$ws.Range("N113").Value = -6728

$ws.Range("H122").Value = 1114.8695
$ws.Range("J122").Value = 4666.3335
$ws.Range("L122").Value = 41997.0015
$ws.Range("N122").Value = -46897.0015


# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2753.5
$ws.Range("I122").Value = 1601.75
$ws.Range("J122").Value = 4135.6
$ws.Range("K122").Value = 4805.25
$ws.Range("L122").Value = 12406.8
$ws.Range("M122").Value = -2355.25
$ws.Range("N122").Value = -17306.8

$ws.Range("H132").Value = 2535.4688
$ws.Range("I132").Value = 2141.15
$ws.Range("J132").Value = 3192.6667
$ws.Range("K132").Value = 6423.450000000001
$ws.Range("L132").Value = 9578.000100000001
$ws.Range("M132").Value = -3893.450000000001
$ws.Range("N132").Value = -14638.0001

$ws.Range("H135").Value = 41146.668
$ws.Range("J135").Value = 41146.668
$ws.Range("L135").Value = 41146.668
$ws.Range("N135").Value = -51286.668


# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 168000
$ws.Range("I40").Value = 168000
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 168000
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -167864

$ws.Range("H46").Value = 2333.9285
$ws.Range("I46").Value = 1834.375
$ws.Range("K46").Value = 1834.375
$ws.Range("M46").Value = -1646.375

$ws.Range("H82").Value = 1615.9354
$ws.Range("I82").Value = 1463.3182
$ws.Range("J82").Value = 1989
$ws.Range("K82").Value = 1463.3182
$ws.Range("L82").Value = 1989
$ws.Range("M82").Value = -1102.3182
$ws.Range("N82").Value = -2711

$ws.Range("H85").Value = 1615.9354
$ws.Range("I85").Value = 1463.3182
$ws.Range("J85").Value = 1989
$ws.Range("K85").Value = 1463.3182
$ws.Range("L85").Value = 1989
$ws.Range("M85").Value = -215.3181999999999
$ws.Range("N85").Value = -4485

$ws.Range("H127").Value = 50736.875
$ws.Range("J127").Value = 50736.875
$ws.Range("L127").Value = 50736.875
$ws.Range("N127").Value = -60656.875

$ws.Range("H132").Value = 6973.3906
$ws.Range("I132").Value = 7294.744
$ws.Range("J132").Value = 6315.381
$ws.Range("K132").Value = 21884.232
$ws.Range("L132").Value = 18946.143
$ws.Range("M132").Value = -19354.232
$ws.Range("N132").Value = -24006.143

$ws.Range("N40").ClearContents()

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9983
$ws.Range("J15").Value = 9983
$ws.Range("L15").Value = 9983
$ws.Range("N15").Value = -10559

$ws.Range("H54").Value = 30590
$ws.Range("J54").Value = 30590
$ws.Range("L54").Value = 30590
$ws.Range("N54").Value = -31630

$ws.Range("H81").Value = 2153
$ws.Range("I81").Value = 1300.1666
$ws.Range("J81").Value = 2618.182
$ws.Range("K81").Value = 2600.3332
$ws.Range("L81").Value = 5236.364
$ws.Range("M81").Value = -1539.3332
$ws.Range("N81").Value = -7358.364

$ws.Range("H84").Value = 2153
$ws.Range("I84").Value = 1300.1666
$ws.Range("J84").Value = 2618.182
$ws.Range("K84").Value = 13001.666
$ws.Range("L84").Value = 26181.82
$ws.Range("M84").Value = -7697.666000000001
$ws.Range("N84").Value = -36789.82

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0

$ws.Range("H100").Value = 828.75
$ws.Range("I100").Value = 818.5714
$ws.Range("J100").Value = 900
$ws.Range("K100").Value = 1637.1428
$ws.Range("L100").Value = 1800
$ws.Range("M100").Value = -1096.1428
$ws.Range("N100").Value = -2882

$ws.Range("H103").Value = 30602
$ws.Range("J103").Value = 30602
$ws.Range("L103").Value = 30602
$ws.Range("N103").Value = -32946

$ws.Range("H132").Value = 2165.0815
$ws.Range("I132").Value = 1844.8966
$ws.Range("J132").Value = 2629.35
$ws.Range("K132").Value = 5534.6898
$ws.Range("L132").Value = 7888.049999999999
$ws.Range("M132").Value = -3004.6898
$ws.Range("N132").Value = -12948.05

$ws.Range("N97").ClearContents()
